$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "2026-01-28 09:59"
$ws.Range("B6").Value = 39
$ws.Range("C6").Value = 7
